# Mask sensitive data with Base64 encoding, per commit message
# "Use Base64 encryption to mask sensitive data"

$wb = $excel.ActiveWorkbook

# --- 1) RunManager sheet: oauth2CC client credentials string (C2) -> base64 parts ---
$runSheet = $wb.Worksheets.Item("RunManager")
$newOAuth = "client_id:cmVzdF9hc3N1cmVkX29hdXRoX2RlbW9fYXBw;client_secret:MDVjM2U3YjI4OGUwZjljZDJhNTFhZmExZGM2NjBkMTE=;grant_type:Y2xpZW50X2NyZWRlbnRpYWxz"
$runSheet.Range("C2").Value = $newOAuth

# --- 2) RunManager sheet: row 4 ("BookTests") - move "basic" auth value from
#        authType (B4) to headers (E4), encoded in base64 ---
$newBasic = "QmFzaWMgYldWeVgyZHlaWGs2VTJobGNHaGxjbVJBTVRJeg=="

$runSheet.Range("B4").ClearContents()
$runSheet.Range("E4").Value = $newBasic

# Copy the direct cell formatting (borders/fill/numfmt) from C4, which already
# carries the same visual style that E4 should end up with.
$runSheet.Range("C4").Copy()
$runSheet.Range("E4").PasteSpecial(-4122)  # xlPasteFormats
